# Front polish 20 mins
# Applies:
#  - D14: "15m" -> "1h"
#  - new values for D18, D19, D20, D22, D23, D30 (realno utroseno vreme column)
#  - clears E24 ("transaction type") and E28 ("balance") which also drops
#    those two now-unused shared strings
#  - appends a new "5m" value (implicitly, by being used in D22/D23)
#  - scrolls/selects so the view shows D31 with A13 pinned at the top

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Value = "1h"

$ws.Range("D18").Value = "1h"
$ws.Range("D19").Value = "30m"
$ws.Range("D20").Value = "20m"

$ws.Range("D22").Value = "5m"
$ws.Range("D23").Value = "5m"

$ws.Range("E24").ClearContents()
$ws.Range("E28").ClearContents()

$ws.Range("D30").Value = "20m"

$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D31").Select()
